$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_9_8_2"
$ws.Range("B2").Value = 0.2223819225481615
$ws.Range("C2").Value = 0.06626570326100178
$ws.Range("D2").Value = 0.4048390793495529
$ws.Range("E2").Value = 0.2789542991640961
$ws.Range("F2").Value = 0.8605940937995911
$ws.Range("G2").Value = 0.7609127759933472
$ws.Range("H2").Value = 0.8247852325439453
$ws.Range("I2").Value = 0.7909711599349976

$ws.Range("A3").Value = "model_9_8_4"
$ws.Range("B3").Value = 0.2316744789872431
$ws.Range("C3").Value = 0.008171567671399038
$ws.Range("D3").Value = 0.2988571406073885
$ws.Range("E3").Value = 0.1931022960817259
$ws.Range("F3").Value = 0.8503099083900452
$ws.Range("G3").Value = 0.8082544803619385
$ws.Range("H3").Value = 0.9716569781303406
$ws.Range("I3").Value = 0.8851489424705505

$ws.Range("A4").Value = "model_9_8_3"
$ws.Range("B4").Value = 0.2359974050771871
$ws.Range("C4").Value = 0.08090093816264099
$ws.Range("D4").Value = 0.3829755291995676
$ws.Range("E4").Value = 0.2717126938056112
$ws.Range("F4").Value = 0.8455257415771484
$ws.Range("G4").Value = 0.7489863634109497
$ws.Range("H4").Value = 0.8550841808319092
$ws.Range("I4").Value = 0.7989150881767273

$ws.Range("A5").Value = "model_9_8_5"
$ws.Range("B5").Value = 0.2546399374424804
$ws.Range("C5").Value = -0.04414961281853746
$ws.Range("D5").Value = 0.230951419973387
$ws.Range("E5").Value = 0.1321546639209354
$ws.Range("F5").Value = 0.8248939514160156
$ws.Range("G5").Value = 0.8508917689323425
$ws.Range("H5").Value = 1.065762042999268
$ws.Range("I5").Value = 0.9520071744918823

$ws.Range("A6").Value = "model_9_8_0"
$ws.Range("B6").Value = 0.2551232231753826
$ws.Range("C6").Value = 0.5568190210572888
$ws.Range("D6").Value = 0.6115602197430489
$ws.Range("E6").Value = 0.5947742636032562
$ws.Range("F6").Value = 0.8243590593338013
$ws.Range("G6").Value = 0.3611541986465454
$ws.Range("H6").Value = 0.5383071899414062
$ws.Range("I6").Value = 0.4445236325263977

$ws.Range("A7").Value = "model_9_8_6"
$ws.Range("B7").Value = 0.2568461061922943
$ws.Range("C7").Value = 0.01146742881494567
$ws.Range("D7").Value = 0.2026694970289793
$ws.Range("E7").Value = 0.1372130529640498
$ws.Range("F7").Value = 0.8224524259567261
$ws.Range("G7").Value = 0.8055685758590698
$ws.Range("H7").Value = 1.104955673217773
$ws.Range("I7").Value = 0.9464582800865173

$ws.Range("A8").Value = "model_9_8_1"
$ws.Range("B8").Value = 0.257867391926151
$ws.Range("C8").Value = 0.5115005466103166
$ws.Range("D8").Value = 0.5314491638340335
$ws.Range("E8").Value = 0.5293284079684425
$ws.Range("F8").Value = 0.8213220238685608
$ws.Range("G8").Value = 0.3980848491191864
$ws.Range("H8").Value = 0.6493266224861145
$ws.Range("I8").Value = 0.5163163542747498

$ws.Range("A9").Value = "model_9_8_7"
$ws.Range("B9").Value = 0.2752153649436345
$ws.Range("C9").Value = -0.8713669290705961
$ws.Range("D9").Value = 0.1696217495939353
$ws.Range("E9").Value = -0.2296343802308165
$ws.Range("F9").Value = 0.8021230697631836
$ws.Range("G9").Value = 1.525002479553223
$ws.Range("H9").Value = 1.150753855705261
$ws.Range("I9").Value = 1.348881721496582

$ws.Range("A10").Value = "model_9_8_8"
$ws.Range("B10").Value = 0.3029353102911027
$ws.Range("C10").Value = -0.8928619443649282
$ws.Range("D10").Value = 0.1843094919213473
$ws.Range("E10").Value = -0.2293592458401568
$ws.Range("F10").Value = 0.7714451551437378
$ws.Range("G10").Value = 1.542518973350525
$ws.Range("H10").Value = 1.130399346351624
$ws.Range("I10").Value = 1.348580002784729

$ws.Range("A11").Value = "model_9_8_9"
$ws.Range("B11").Value = 0.3653874855614727
$ws.Range("C11").Value = -0.8150426295014475
$ws.Range("D11").Value = 0.1862756799056046
$ws.Range("E11").Value = -0.1975848654151884
$ws.Range("F11").Value = 0.7023290991783142
$ws.Range("G11").Value = 1.479102849960327
$ws.Range("H11").Value = 1.127674460411072
$ws.Range("I11").Value = 1.313724040985107

$ws.Range("A12").Value = "model_9_8_11"
$ws.Range("B12").Value = 0.5649386743985904
$ws.Range("C12").Value = -0.225660683666526
$ws.Range("D12").Value = 0.2344611715151383
$ws.Range("E12").Value = 0.06285355094172307
$ws.Range("F12").Value = 0.4814846813678741
$ws.Range("G12").Value = 0.9988076090812683
$ws.Range("H12").Value = 1.060898184776306
$ws.Range("I12").Value = 1.028028964996338

$ws.Range("A13").Value = "model_9_8_10"
$ws.Range("B13").Value = 0.5736627276849307
$ws.Range("C13").Value = -0.2389452983427303
$ws.Range("D13").Value = 0.3128675461929177
$ws.Range("E13").Value = 0.1042427120320183
$ws.Range("F13").Value = 0.4718296825885773
$ws.Range("G13").Value = 1.009633421897888
$ws.Range("H13").Value = 0.9522411227226257
$ws.Range("I13").Value = 0.9826259613037109

$ws.Range("A14").Value = "model_9_8_23"
$ws.Range("B14").Value = 0.5847975310803115
$ws.Range("C14").Value = 0.4503647073527868
$ws.Range("D14").Value = -0.8275392540677455
$ws.Range("E14").Value = -0.3026169685436364
$ws.Range("F14").Value = 0.4595067501068115
$ws.Range("G14").Value = 0.4479052722454071
$ws.Range("H14").Value = 2.532638549804688
$ws.Range("I14").Value = 1.428942084312439

$ws.Range("A15").Value = "model_9_8_19"
$ws.Range("B15").Value = 0.6526736761829517
$ws.Range("C15").Value = 0.3836275347718785
$ws.Range("D15").Value = -0.3556186036833235
$ws.Range("E15").Value = -0.04832030408073118
$ws.Range("F15").Value = 0.3843879103660583
$ws.Range("G15").Value = 0.5022903084754944
$ws.Range("H15").Value = 1.878641843795776
$ws.Range("I15").Value = 1.149984240531921

$ws.Range("A16").Value = "model_9_8_18"
$ws.Range("B16").Value = 0.6555001223448682
$ws.Range("C16").Value = 0.3672138702737662
$ws.Range("D16").Value = -0.309691421605111
$ws.Range("E16").Value = -0.02747932256521346
$ws.Range("F16").Value = 0.3812598586082458
$ws.Range("G16").Value = 0.5156659483909607
$ws.Range("H16").Value = 1.814995050430298
$ws.Range("I16").Value = 1.127122044563293

$ws.Range("A17").Value = "model_9_8_24"
$ws.Range("B17").Value = 0.6728885420075468
$ws.Range("C17").Value = 0.4257184162268742
$ws.Range("D17").Value = -0.3293255490828213
$ws.Range("E17").Value = -0.01613690691464242
$ws.Range("F17").Value = 0.362015962600708
$ws.Range("G17").Value = 0.4679898619651794
$ws.Range("H17").Value = 1.842204332351685
$ws.Range("I17").Value = 1.11467969417572

$ws.Range("A18").Value = "model_9_8_12"
$ws.Range("B18").Value = 0.6746146984210148
$ws.Range("C18").Value = 0.1058834812189207
$ws.Range("D18").Value = 0.1323931440009815
$ws.Range("E18").Value = 0.1325660246179416
$ws.Range("F18").Value = 0.3601056635379791
$ws.Range("G18").Value = 0.7286277413368225
$ws.Range("H18").Value = 1.202345967292786
$ws.Range("I18").Value = 0.9515559077262878

$ws.Range("A19").Value = "model_9_8_21"
$ws.Range("B19").Value = 0.6794223017433868
$ws.Range("C19").Value = 0.4227572138375939
$ws.Range("D19").Value = -0.2606515731092771
$ws.Range("E19").Value = 0.02352830205014089
$ws.Range("F19").Value = 0.3547850251197815
$ws.Range("G19").Value = 0.4704029858112335
$ws.Range("H19").Value = 1.747034788131714
$ws.Range("I19").Value = 1.071167945861816

$ws.Range("A20").Value = "model_9_8_20"
$ws.Range("B20").Value = 0.6822406578433018
$ws.Range("C20").Value = 0.4219520674385736
$ws.Range("D20").Value = -0.2406301958635384
$ws.Range("E20").Value = 0.03511591042739626
$ws.Range("F20").Value = 0.3516659438610077
$ws.Range("G20").Value = 0.471059113740921
$ws.Range("H20").Value = 1.71928882598877
$ws.Range("I20").Value = 1.058456420898438

$ws.Range("A21").Value = "model_9_8_17"
$ws.Range("B21").Value = 0.6886499824726686
$ws.Range("C21").Value = 0.2989603520401164
$ws.Range("D21").Value = -0.05980025738479067
$ws.Range("E21").Value = 0.09424255897658806
$ws.Range("F21").Value = 0.3445727527141571
$ws.Range("G21").Value = 0.5712867379188538
$ws.Range("H21").Value = 1.468691229820251
$ws.Range("I21").Value = 0.9935958981513977

$ws.Range("A22").Value = "model_9_8_14"
$ws.Range("B22").Value = 0.6913227870086531
$ws.Range("C22").Value = 0.06956287927527416
$ws.Range("D22").Value = 0.1840857843291588
$ws.Range("E22").Value = 0.1490162570619028
$ws.Range("F22").Value = 0.3416147232055664
$ws.Range("G22").Value = 0.7582257986068726
$ws.Range("H22").Value = 1.130709409713745
$ws.Range("I22").Value = 0.9335103034973145

$ws.Range("A23").Value = "model_9_8_15"
$ws.Range("B23").Value = 0.6915562255113754
$ws.Range("C23").Value = 0.1300518129848786
$ws.Range("D23").Value = 0.09947248755211213
$ws.Range("E23").Value = 0.122504295443382
$ws.Range("F23").Value = 0.3413564264774323
$ws.Range("G23").Value = 0.708932638168335
$ws.Range("H23").Value = 1.247968077659607
$ws.Range("I23").Value = 0.9625933170318604

$ws.Range("A24").Value = "model_9_8_13"
$ws.Range("B24").Value = 0.6942178477741929
$ws.Range("C24").Value = 0.07468328536297153
$ws.Range("D24").Value = 0.2092986151467112
$ws.Range("E24").Value = 0.1660177103374717
$ws.Range("F24").Value = 0.3384107649326324
$ws.Range("G24").Value = 0.7540531754493713
$ws.Range("H24").Value = 1.095768928527832
$ws.Range("I24").Value = 0.9148601293563843

$ws.Range("A25").Value = "model_9_8_16"
$ws.Range("B25").Value = 0.695906255125834
$ws.Range("C25").Value = 0.1990147259565227
$ws.Range("D25").Value = 0.07477401466247657
$ws.Range("E25").Value = 0.1349403963021975
$ws.Range("F25").Value = 0.3365421891212463
$ws.Range("G25").Value = 0.6527338027954102
$ws.Range("H25").Value = 1.282195568084717
$ws.Range("I25").Value = 0.948951244354248

$ws.Range("A26").Value = "model_9_8_22"
$ws.Range("B26").Value = 0.6962350887706158
$ws.Range("C26").Value = 0.4298422489783221
$ws.Range("D26").Value = -0.1811771395446669
$ws.Range("E26").Value = 0.07356141706243347
$ws.Range("F26").Value = 0.3361782431602478
$ws.Range("G26").Value = 0.4646292924880981
$ws.Range("H26").Value = 1.636897563934326
$ws.Range("I26").Value = 1.016282558441162
